$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text number format on the cells we touch so Excel does not
# auto-convert the numeric-looking / percent-looking strings into
# actual numbers (these columns store plain text values).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '309.06'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '2.07%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '38.71'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '8.48%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.092'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '1.15%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08179'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '3.24%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.978'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '6.33%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '7.886'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '1.77%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9318'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '1.39%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1401'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '3.80%'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '3.75%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.09215'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '2.18%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03442'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '0.22%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09848'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '0.44%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001410'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.84%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.006116'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.56%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.667'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-1.78%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.183'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '1.97%'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '5.71%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3450'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.25%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1330'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '2.32%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.812'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-6.53%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2454'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '2.99%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04469'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '1.63%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001239'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.43%'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-9.52%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001302'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '0.36%'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '9.86%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05170'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-1.37%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007474'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-1.70%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.01000'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-0.92%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1369'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '1.82%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002133'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-0.57%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.009674'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-4.62%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006317'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '2.95%'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.34%'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '1.94%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.001602'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-3.26%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002103'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.34%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002003'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.34%'
